$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.40311104990069
$ws.Range("C2").Value = 10.34457366267983
$ws.Range("D2").Value = 5.997725142654222
$ws.Range("E2").Value = 12.53742077145915
$ws.Range("G2").Value = 34.50002494905634
$ws.Range("H2").Value = 15.78068430717363
$ws.Range("K2").Value = 8.781178598514009
$ws.Range("L2").Value = 9.568679151578802
$ws.Range("O2").Value = 24.77787077309823
$ws.Range("B3").Value = 12.13252846446118
$ws.Range("C3").Value = 10.36116578130155
$ws.Range("D3").Value = 5.88102639286096
$ws.Range("E3").Value = 12.55246835435037
$ws.Range("G3").Value = 34.57191597252015
$ws.Range("H3").Value = 15.8295670183605
$ws.Range("K3").Value = 8.58070110859977
$ws.Range("L3").Value = 9.553648080333989
$ws.Range("O3").Value = 24.85442283197309
$ws.Range("B4").Value = 11.96544408291313
$ws.Range("C4").Value = 10.37206965038587
$ws.Range("D4").Value = 5.80993606952778
$ws.Range("E4").Value = 12.56408971423606
$ws.Range("G4").Value = 34.6261391637868
$ws.Range("H4").Value = 15.86206819510987
$ws.Range("K4").Value = 8.456366744499906
$ws.Range("L4").Value = 9.546044627345621
$ws.Range("O4").Value = 24.90653032640696
$ws.Range("B5").Value = 11.89721607905655
$ws.Range("C5").Value = 10.37669366622602
$ws.Range("D5").Value = 5.781152508491803
$ws.Range("E5").Value = 12.56942438909771
$ws.Range("G5").Value = 34.65076190972064
$ws.Range("H5").Value = 15.8759379212868
$ws.Range("K5").Value = 8.405458502425397
$ws.Range("L5").Value = 9.54335747397489
$ws.Range("O5").Value = 24.92904564855782
$ws.Range("B6").Value = 11.88588135370429
$ws.Range("C6").Value = 10.3774724030104
$ws.Range("D6").Value = 5.776385587079695
$ws.Range("E6").Value = 12.57034637180988
$ws.Range("G6").Value = 34.65500280714215
$ws.Range("H6").Value = 15.87827873707642
$ws.Range("K6").Value = 8.396992800112816
$ws.Range("L6").Value = 9.542936184925194
$ws.Range("O6").Value = 24.93286160344179
$ws.Range("B7").Value = 11.96452437032321
$ws.Range("C7").Value = 10.37213127957435
$ws.Range("D7").Value = 5.809547070588571
$ws.Range("E7").Value = 12.56415923497182
$ws.Range("G7").Value = 34.62646101802039
$ws.Range("H7").Value = 15.86225271574305
$ws.Range("K7").Value = 8.455681058097992
$ws.Range("L7").Value = 9.546006718925176
$ws.Range("O7").Value = 24.90682879207722
$ws.Range("B8").Value = 12.31006823114769
$ws.Range("C8").Value = 10.35014626648762
$ws.Range("D8").Value = 5.957397202134868
$ws.Range("E8").Value = 12.5421147517408
$ws.Range("G8").Value = 34.52271549625559
$ws.Range("H8").Value = 15.79702280747877
$ws.Range("K8").Value = 8.712353053878559
$ws.Range("L8").Value = 9.563160621720712
$ws.Range("O8").Value = 24.80320512767066
$ws.Range("B9").Value = 12.97604591311931
$ws.Range("C9").Value = 10.31269479128342
$ws.Range("D9").Value = 6.249866277704904
$ws.Range("E9").Value = 12.51778614964307
$ws.Range("G9").Value = 34.39961566399957
$ws.Range("H9").Value = 15.68884653378466
$ws.Range("K9").Value = 9.202882808536298
$ws.Range("L9").Value = 9.609576502324007
$ws.Range("O9").Value = 24.64060309380287
$ws.Range("B10").Value = 13.45300941767567
$ws.Range("C10").Value = 10.28860048138828
$ws.Range("D10").Value = 6.463773808295074
$ws.Range("E10").Value = 12.51142578124281
$ws.Range("G10").Value = 34.35856799065531
$ws.Range("H10").Value = 15.62141127844662
$ws.Range("K10").Value = 9.551747594218989
$ws.Range("L10").Value = 9.651288066867053
$ws.Range("O10").Value = 24.54602931128908
$ws.Range("B11").Value = 13.66630105430167
$ws.Range("C11").Value = 10.27837600439951
$ws.Range("D11").Value = 6.560372402325252
$ws.Range("E11").Value = 12.51102700117554
$ws.Range("G11").Value = 34.35068458843995
$ws.Range("H11").Value = 15.59334975016935
$ws.Range("K11").Value = 9.707240261938058
$ws.Range("L11").Value = 9.671872275580208
$ws.Range("O11").Value = 24.50843686152782
$ws.Range("B12").Value = 13.74646217684199
$ws.Range("C12").Value = 10.27460963328301
$ws.Range("D12").Value = 6.596810477016862
$ws.Range("E12").Value = 12.51123390267959
$ws.Range("G12").Value = 34.34925411717037
$ws.Range("H12").Value = 15.58309980610407
$ws.Range("K12").Value = 9.765606109269859
$ws.Range("L12").Value = 9.679894247348498
$ws.Range("O12").Value = 24.49498437590276
$ws.Range("B13").Value = 13.72922626701285
$ws.Range("C13").Value = 10.27541610722201
$ws.Range("D13").Value = 6.588969798668
$ws.Range("E13").Value = 12.5111734395106
$ws.Range("G13").Value = 34.34949300562804
$ws.Range("H13").Value = 15.58529057212226
$ws.Range("K13").Value = 9.753059748002384
$ws.Range("L13").Value = 9.678156538816028
$ws.Range("O13").Value = 24.49784675385444
$ws.Range("B14").Value = 13.67290856888173
$ws.Range("C14").Value = 10.27806403258197
$ws.Range("D14").Value = 6.563373244949503
$ws.Range("E14").Value = 12.51103685534408
$ws.Range("G14").Value = 34.35053572781788
$ws.Range("H14").Value = 15.5924989368541
$ws.Range("K14").Value = 9.712052691868244
$ws.Range("L14").Value = 9.672527720300227
$ws.Range("O14").Value = 24.50731441378681
$ws.Range("B15").Value = 13.63833096621232
$ws.Range("C15").Value = 10.27969967867598
$ws.Range("D15").Value = 6.547675017706851
$ws.Range("E15").Value = 12.51099977683875
$ws.Range("G15").Value = 34.35137698245265
$ws.Range("H15").Value = 15.59696328868477
$ws.Range("K15").Value = 9.686865956680561
$ws.Range("L15").Value = 9.669109361184072
$ws.Range("O15").Value = 24.51321565556878
$ws.Range("B16").Value = 13.43898940135637
$ws.Range("C16").Value = 10.28928345035972
$ws.Range("D16").Value = 6.457443148283537
$ws.Range("E16").Value = 12.51150197534791
$ws.Range("G16").Value = 34.35930061102295
$ws.Range("H16").Value = 15.62329781365024
$ws.Range("K16").Value = 9.541516475780773
$ws.Range("L16").Value = 9.649974868275338
$ws.Range("O16").Value = 24.548595514519
$ws.Range("B17").Value = 13.31569939614084
$ws.Range("C17").Value = 10.29535101324365
$ws.Range("D17").Value = 6.401878262425812
$ws.Range("E17").Value = 12.51244846861513
$ws.Range("G17").Value = 34.36692767582377
$ws.Range("H17").Value = 15.64012319830991
$ws.Range("K17").Value = 9.451487570685382
$ws.Range("L17").Value = 9.638645742289684
$ws.Range("O17").Value = 24.57169217273414
$ws.Range("B18").Value = 13.24444407711929
$ws.Range("C18").Value = 10.29891023039499
$ws.Range("D18").Value = 6.369854322233357
$ws.Range("E18").Value = 12.51322772923622
$ws.Range("G18").Value = 34.37232999293192
$ws.Range("H18").Value = 15.65004684344079
$ws.Range("K18").Value = 9.399406408155892
$ws.Range("L18").Value = 9.632281285139033
$ws.Range("O18").Value = 24.58548770129541
$ws.Range("B19").Value = 13.22026196036693
$ws.Range("C19").Value = 10.30012723867576
$ws.Range("D19").Value = 6.359001736823974
$ws.Range("E19").Value = 12.51353193276451
$ws.Range("G19").Value = 34.37433340450682
$ws.Range("H19").Value = 15.65344908284062
$ws.Range("K19").Value = 9.381723012671431
$ws.Range("L19").Value = 9.630152578917221
$ws.Range("O19").Value = 24.59024631604028
$ws.Range("B20").Value = 13.32885982602963
$ws.Range("C20").Value = 10.2946979401333
$ws.Range("D20").Value = 6.40780020840907
$ws.Range("E20").Value = 12.51232341060604
$ws.Range("G20").Value = 34.36601064043695
$ws.Range("H20").Value = 15.63830663177591
$ws.Range("K20").Value = 9.461102646628261
$ws.Range("L20").Value = 9.639836070847316
$ws.Range("O20").Value = 24.56918059719255
$ws.Range("B21").Value = 13.68946753511942
$ws.Range("C21").Value = 10.27728341522763
$ws.Range("D21").Value = 6.570895725574923
$ws.Range("E21").Value = 12.51106726723746
$ws.Range("G21").Value = 34.35018723821455
$ws.Range("H21").Value = 15.5903714490925
$ws.Range("K21").Value = 9.72411187306832
$ws.Range("L21").Value = 9.674174910347702
$ws.Range("O21").Value = 24.50451226660422
$ws.Range("B22").Value = 13.92157032430205
$ws.Range("C22").Value = 10.2665162636509
$ws.Range("D22").Value = 6.676645146087877
$ws.Range("E22").Value = 12.51233196933639
$ws.Range("G22").Value = 34.34890892505648
$ws.Range("H22").Value = 15.56123681105617
$ws.Range("K22").Value = 9.892974394606169
$ws.Range("L22").Value = 9.69793934581257
$ws.Range("O22").Value = 24.46681236725093
$ws.Range("B23").Value = 13.79804432278857
$ws.Range("C23").Value = 10.27220683549483
$ws.Range("D23").Value = 6.620294391253536
$ws.Range("E23").Value = 12.511466453277
$ws.Range("G23").Value = 34.34876114462706
$ws.Range("H23").Value = 15.5765856886858
$ws.Range("K23").Value = 9.8031434391415
$ws.Range("L23").Value = 9.685136325510525
$ws.Range("O23").Value = 24.48651516950767
$ws.Range("B24").Value = 13.32291115537213
$ws.Range("C24").Value = 10.29499297362163
$ws.Range("D24").Value = 6.405123139498799
$ws.Range("E24").Value = 12.5123792169384
$ws.Range("G24").Value = 34.36642206319434
$ws.Range("H24").Value = 15.63912712119905
$ws.Range("K24").Value = 9.456756670638157
$ws.Range("L24").Value = 9.639297459470152
$ws.Range("O24").Value = 24.57031447036749
$ws.Range("B25").Value = 12.79768546617054
$ws.Range("C25").Value = 10.32222345838668
$ws.Range("D25").Value = 6.170739254369307
$ws.Range("E25").Value = 12.52234439457123
$ws.Range("G25").Value = 34.42426880080818
$ws.Range("H25").Value = 15.7159969819904
$ws.Range("K25").Value = 9.071950965901671
$ws.Range("L25").Value = 9.595669314307912
$ws.Range("O25").Value = 24.68022974404089
